$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 4).Value = '64.335.33'
$ws.Cells.Item(2, 5).Value = '  +1.67%  '

$ws.Cells.Item(3, 4).Value = '3.086.71'
$ws.Cells.Item(3, 5).Value = '  +0.89%  '

$ws.Cells.Item(4, 5).Value = '  +0.03%  '

$ws.Cells.Item(5, 4).Formula = '="560.87"'
$ws.Cells.Item(5, 5).Value = '  +2.17%  '

$ws.Cells.Item(6, 4).Formula = '="145.48"'
$ws.Cells.Item(6, 5).Value = '  +4.03%  '

$ws.Cells.Item(7, 4).Formula = '="1.00"'
$ws.Cells.Item(7, 5).Value = '  +0.04%  '

$ws.Cells.Item(8, 4).Value = '3.084.04'
$ws.Cells.Item(8, 5).Value = '  +1.01%  '

$ws.Cells.Item(9, 4).Formula = '="0.506"'
$ws.Cells.Item(9, 5).Value = '  +0.96%  '

$ws.Cells.Item(10, 5).Value = '  +2.25%  '

$ws.Cells.Item(11, 4).Formula = '="6.13"'
$ws.Cells.Item(11, 5).Value = '  -4.32%  '

$ws.Cells.Item(12, 5).Value = '  +4.35%  '

$ws.Cells.Item(13, 5).Value = '  +1.08%  '

$ws.Cells.Item(14, 4).Formula = '="35.25"'
$ws.Cells.Item(14, 5).Value = '  +1.57%  '

$ws.Cells.Item(15, 4).Value = '3.582.02'
$ws.Cells.Item(15, 5).Value = '  +0.88%  '

$ws.Cells.Item(16, 4).Value = '64.363.49'
$ws.Cells.Item(16, 5).Value = '  +1.66%  '

$ws.Cells.Item(17, 4).Value = '3.080.82'
$ws.Cells.Item(17, 5).Value = '  +0.75%  '

$ws.Cells.Item(19, 5).Value = '  +0.67%  '

$ws.Cells.Item(20, 4).Formula = '="480.30"'
$ws.Cells.Item(20, 5).Value = '  -0.22%  '

$ws.Cells.Item(21, 5).Value = '  +2.06%  '

$ws.Cells.Item(22, 4).Formula = '="0.676"'
$ws.Cells.Item(22, 5).Value = '  +0.55%  '

$ws.Cells.Item(23, 5).Value = '  +4.89%  '

$ws.Cells.Item(24, 4).Formula = '="13.77"'
$ws.Cells.Item(24, 5).Value = '  +9.89%  '

$ws.Cells.Item(25, 4).Formula = '="81.31"'
$ws.Cells.Item(25, 5).Value = '  +0.83%  '

$ws.Cells.Item(26, 4).Formula = '="0.996"'
$ws.Cells.Item(26, 5).Value = '  -0.39%  '

$ws.Cells.Item(27, 4).Formula = '="2.82"'
$ws.Cells.Item(27, 5).Value = '  +2.43%  '

$ws.Cells.Item(28, 4).Formula = '="8.04"'
$ws.Cells.Item(28, 5).Value = '  +1.49%  '

$ws.Cells.Item(29, 4).Formula = '="2.08"'
$ws.Cells.Item(29, 5).Value = '  +5.07%  '

$ws.Cells.Item(30, 4).Formula = '="0.999"'
$ws.Cells.Item(30, 5).Value = '  +0.06%  '

$ws.Cells.Item(31, 4).Formula = '="26.26"'
$ws.Cells.Item(31, 5).Value = '  +1.01%  '

$ws.Cells.Item(32, 4).Formula = '="1.15"'
$ws.Cells.Item(32, 5).Value = '  -0.14%  '

$ws.Cells.Item(33, 4).Formula = '="2.50"'
$ws.Cells.Item(33, 5).Value = '  +3.10%  '

$ws.Cells.Item(34, 5).Value = '  -1.70%  '

$ws.Cells.Item(35, 5).Value = '  +4.04%  '

$ws.Cells.Item(36, 4).Formula = '="55.78"'
$ws.Cells.Item(36, 5).Value = '  +0.67%  '

$ws.Cells.Item(37, 4).Formula = '="3.07"'
$ws.Cells.Item(37, 5).Value = '  +18.78%  '

$ws.Cells.Item(38, 4).Formula = '="458.25"'
$ws.Cells.Item(38, 5).Value = '  -1.77%  '

$ws.Cells.Item(39, 5).Value = '  +3.12%  '

$ws.Cells.Item(40, 4).Formula = '="0.0825"'
$ws.Cells.Item(40, 5).Value = '  +0.92%  '

$ws.Cells.Item(41, 4).Value = '2.975.27'
$ws.Cells.Item(41, 5).Value = '  -2.92%  '

$ws.Cells.Item(42, 5).Value = '  +0.23%  '

$ws.Cells.Item(43, 5).Value = '  -2.23%  '

$ws.Cells.Item(44, 4).Formula = '="28.04"'
$ws.Cells.Item(44, 5).Value = '  -0.91%  '

$ws.Cells.Item(45, 5).Value = '  +3.47%  '

$ws.Cells.Item(46, 5).Value = '  +4.98%  '

$ws.Cells.Item(48, 5).Value = '  +2.53%  '

$ws.Cells.Item(49, 4).Formula = '="120.93"'
$ws.Cells.Item(49, 5).Value = '  +3.34%  '

$ws.Cells.Item(50, 5).Value = '  +1.70%  '

$ws.Cells.Item(51, 5).Value = '  +1.11%  '
